# Adds undo and redo implementation
# Nudges a handful of shapes on the UndoAndRedoImplementation stack-diagram
# slides so their yellow "command" boxes / labels sit correctly inside the
# stack outlines instead of slightly overflowing them.

$p = $ppt.ActivePresentation

# EMU -> point helper (PowerPoint COM positions/sizes are expressed in points)
function EmuToPt($emu) { return $emu / 12700.0 }

# --- Slide 1: "RedoStack" label (TextBox 13) -------------------------------
# Move it down so it lines up with the "UndoStack" label on the same slide.
$s1 = $p.Slides.Item(1)
$redoStack1 = $s1.Shapes.Item("TextBox 13")
$redoStack1.Left = EmuToPt 6871939
$redoStack1.Top = EmuToPt 5978030
$redoStack1.Width = EmuToPt 1157240
$redoStack1.Height = EmuToPt 369332

# --- Slide 2: "a:AddExpenseCommand" box (Rectangle 9) -----------------------
# Nudge down slightly and widen a touch so the box/text sit inside the stack.
$s2 = $p.Slides.Item(2)
$rect9_s2 = $s2.Shapes.Item("Rectangle 9")
$rect9_s2.Left = EmuToPt 3972040
$rect9_s2.Top = EmuToPt 5518694
$rect9_s2.Width = EmuToPt 1717446
$rect9_s2.Height = EmuToPt 316727

# --- Slide 3: "a:AddExpenseCommand" box (Rectangle 9) -----------------------
$s3 = $p.Slides.Item(3)
$rect9_s3 = $s3.Shapes.Item("Rectangle 9")
$rect9_s3.Left = EmuToPt 3972040
$rect9_s3.Top = EmuToPt 5522504
$rect9_s3.Width = EmuToPt 1707127
$rect9_s3.Height = EmuToPt 316727

# --- Slide 3: "u:UpdateCommand" box (Rectangle 14) --------------------------
$rect14_s3 = $s3.Shapes.Item("Rectangle 14")
$rect14_s3.Left = EmuToPt 3972039
$rect14_s3.Top = EmuToPt 5205777
$rect14_s3.Width = EmuToPt 1707127
$rect14_s3.Height = EmuToPt 316727
